# Auto-applied update of market-price-derived columns (H-N) across sheets,
# as produced by the scheduled runner that refreshes Leve profit calculations.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("L5").Value2 = 100
$ws.Range("N5").Value2 = -330
$ws.Range("H5").Value2 = 82.75
$ws.Range("J5").Value2 = 100
$ws.Range("K8").Value2 = 42.59999999999999
$ws.Range("M8").Value2 = 96.40000000000001
$ws.Range("I8").Value2 = 14.2
$ws.Range("H8").Value2 = 14.2
$ws.Range("J51").Value2 = 0
$ws.Range("L51").Value2 = 0
$ws.Range("H51").Value2 = 0
$ws.Range("N51").ClearContents()
$ws.Range("I98").Value2 = 2943.3845
$ws.Range("L98").Value2 = 3296.4285
$ws.Range("N98").Value2 = -6292.4285
$ws.Range("H98").Value2 = 3066.95
$ws.Range("K98").Value2 = 2943.3845
$ws.Range("J98").Value2 = 3296.4285
$ws.Range("M98").Value2 = -1445.3845
$ws.Range("H107").Value2 = 3115
$ws.Range("N107").Value2 = -9266.25
$ws.Range("K107").Value2 = 1266
$ws.Range("J107").Value2 = 5426.25
$ws.Range("M107").Value2 = 654
$ws.Range("L107").Value2 = 5426.25
$ws.Range("I107").Value2 = 1266
$ws.Range("J116").Value2 = 6483
$ws.Range("I116").Value2 = 5722.6665
$ws.Range("L116").Value2 = 6483
$ws.Range("N116").Value2 = -13367
$ws.Range("H116").Value2 = 6157.143
$ws.Range("K116").Value2 = 5722.6665
$ws.Range("M116").Value2 = -2280.6665
$ws.Range("H122").Value2 = 3066.95
$ws.Range("N122").Value2 = -14789.2855
$ws.Range("K122").Value2 = 8830.1535
$ws.Range("M122").Value2 = -6380.1535
$ws.Range("J122").Value2 = 3296.4285
$ws.Range("I122").Value2 = 2943.3845
$ws.Range("L122").Value2 = 9889.2855
$ws.Range("I141").Value2 = 5447.75
$ws.Range("H141").Value2 = 5447.75
$ws.Range("M141").Value2 = -11163.25
$ws.Range("K141").Value2 = 16343.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M88").Value2 = -561
$ws.Range("J88").Value2 = 5629.3076
$ws.Range("I88").Value2 = 967
$ws.Range("L88").Value2 = 5629.3076
$ws.Range("H88").Value2 = 4532.294
$ws.Range("N88").Value2 = -6441.3076
$ws.Range("K88").Value2 = 967
$ws.Range("L91").Value2 = 5629.3076
$ws.Range("I91").Value2 = 967
$ws.Range("H91").Value2 = 4532.294
$ws.Range("N91").Value2 = -8437.3076
$ws.Range("K91").Value2 = 967
$ws.Range("J91").Value2 = 5629.3076
$ws.Range("M91").Value2 = 437

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J11").Value2 = 499
$ws.Range("L11").Value2 = 499
$ws.Range("H11").Value2 = 267.66666
$ws.Range("N11").Value2 = -779

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value2 = 4005.5
$ws.Range("H122").Value2 = 1187.9231
$ws.Range("N122").Value2 = -9467.5
$ws.Range("K122").Value2 = 3117.6666
$ws.Range("M122").Value2 = -667.6665999999996
$ws.Range("J122").Value2 = 1522.5
$ws.Range("I122").Value2 = 1039.2222
$ws.Range("L122").Value2 = 4567.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L10").Value2 = 11994
$ws.Range("I10").Value2 = 19.928572
$ws.Range("H10").Value2 = 285.13333
$ws.Range("N10").Value2 = -12272
$ws.Range("K10").Value2 = 59.78571599999999
$ws.Range("J10").Value2 = 3998
$ws.Range("M10").Value2 = 79.21428400000001
$ws.Range("J37").Value2 = 97518
$ws.Range("L37").Value2 = 292554
$ws.Range("N37").Value2 = -292778
$ws.Range("H37").Value2 = 97518
$ws.Range("K50").Value2 = 1647.5001
$ws.Range("M50").Value2 = -1166.5001
$ws.Range("I50").Value2 = 549.1667
$ws.Range("H50").Value2 = 549.1667
$ws.Range("H53").Value2 = 549.1667
$ws.Range("K53").Value2 = 1647.5001
$ws.Range("M53").Value2 = -1166.5001
$ws.Range("I53").Value2 = 549.1667
$ws.Range("N68").Value2 = -3916.6667
$ws.Range("H68").Value2 = 957.1667
$ws.Range("J68").Value2 = 764.8889
$ws.Range("L68").Value2 = 2294.6667
$ws.Range("H71").Value2 = 957.1667
$ws.Range("N71").Value2 = -14996.0001
$ws.Range("J71").Value2 = 764.8889
$ws.Range("L71").Value2 = 6884.0001
$ws.Range("L92").Value2 = 7500
$ws.Range("I92").Value2 = 2500
$ws.Range("H92").Value2 = 2500
$ws.Range("N92").Value2 = -9996
$ws.Range("K92").Value2 = 7500
$ws.Range("J92").Value2 = 2500
$ws.Range("M92").Value2 = -6252
$ws.Range("H97").Value2 = 903.36365
$ws.Range("N97").Value2 = -2403.28568
$ws.Range("J97").Value2 = 470.42856
$ws.Range("L97").Value2 = 1411.28568
$ws.Range("I98").Value2 = 0
$ws.Range("H98").Value2 = 0
$ws.Range("K98").Value2 = 0
$ws.Range("M98").ClearContents()
$ws.Range("H113").Value2 = 1056
$ws.Range("N113").Value2 = -7662.5
$ws.Range("J113").Value2 = 1107.5
$ws.Range("L113").Value2 = 3322.5
$ws.Range("H122").Value2 = 1396.6666
$ws.Range("K122").Value2 = 9900
$ws.Range("M122").Value2 = -7450
$ws.Range("I122").Value2 = 1100

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value2 = 3961.6667
$ws.Range("N97").ClearContents()
$ws.Range("K97").Value2 = 3961.6667
$ws.Range("M97").Value2 = -3465.6667
$ws.Range("J97").Value2 = 0
$ws.Range("I97").Value2 = 3961.6667
$ws.Range("L97").Value2 = 0
$ws.Range("H122").Value2 = 15628822
$ws.Range("N122").Value2 = -22402
$ws.Range("K122").Value2 = 62509452
$ws.Range("M122").Value2 = -62507002
$ws.Range("J122").Value2 = 5834
$ws.Range("I122").Value2 = 20836484
$ws.Range("L122").Value2 = 17502
$ws.Range("J131").Value2 = 0
$ws.Range("L131").Value2 = 0
$ws.Range("H131").Value2 = 0
$ws.Range("N131").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value2 = 3332.6667
$ws.Range("N12").Value2 = -3672.6667
$ws.Range("J12").Value2 = 3332.6667
$ws.Range("L12").Value2 = 3332.6667
$ws.Range("K46").Value2 = 2958.3333
$ws.Range("J46").Value2 = 3799.8
$ws.Range("M46").Value2 = -2770.3333
$ws.Range("I46").Value2 = 2958.3333
$ws.Range("L46").Value2 = 3799.8
$ws.Range("N46").Value2 = -4175.8
$ws.Range("H46").Value2 = 3484.25
$ws.Range("N55").Value2 = -3757
$ws.Range("H55").Value2 = 2269
$ws.Range("J55").Value2 = 3411
$ws.Range("L55").Value2 = 3411
$ws.Range("H132").Value2 = 5186.129
$ws.Range("K132").Value2 = 14546.8932
$ws.Range("M132").Value2 = -12016.8932
$ws.Range("I132").Value2 = 4848.9644

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J121").Value2 = 30000
$ws.Range("L121").Value2 = 30000
$ws.Range("H121").Value2 = 30000
$ws.Range("N121").Value2 = -33494
$ws.Range("H132").Value2 = 1001.38464
$ws.Range("K132").Value2 = 2868.5454
$ws.Range("M132").Value2 = -338.5454
$ws.Range("I132").Value2 = 956.1818
